$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3
$ws.Range("F5").Value = 837
$ws.Range("F6").Value = 511
$ws.Range("F7").Value = 282
$ws.Range("F8").Value = 9505
$ws.Range("F11").Value = 692
$ws.Range("F12").Value = 1973
$ws.Range("F14").Value = 962
$ws.Range("F15").Value = 2615
$ws.Range("F16").Value = 129
$ws.Range("F17").Value = 3926
$ws.Range("F18").Value = 311
$ws.Range("F19").Value = 141
$ws.Range("F20").Value = 126
$ws.Range("F21").Value = 205
$ws.Range("F22").Value = 236
$ws.Range("F23").Value = 20
$ws.Range("F25").Value = 71
$ws.Range("F27").Value = 564
$ws.Range("F29").Value = 2144
$ws.Range("F32").Value = 467
$ws.Range("F35").Value = 180
$ws.Range("F36").Value = 345
$ws.Range("F37").Value = 161

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 19
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 21

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 191
$ws.Range("F3").Value = 976

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 191
$ws.Range("F4").Value = 976
$ws.Range("F5").Value = 3
$ws.Range("F8").Value = 837
$ws.Range("F9").Value = 511
$ws.Range("F10").Value = 282
$ws.Range("F11").Value = 9505
$ws.Range("F14").Value = 692
$ws.Range("F15").Value = 1974
$ws.Range("F17").Value = 962
$ws.Range("F19").Value = 2615
$ws.Range("F20").Value = 129
$ws.Range("F21").Value = 3926
$ws.Range("F22").Value = 311
$ws.Range("F23").Value = 141
$ws.Range("F24").Value = 126
$ws.Range("F25").Value = 205
$ws.Range("F26").Value = 236
$ws.Range("F27").Value = 20
$ws.Range("F28").Value = 19
$ws.Range("F30").Value = 71
$ws.Range("F32").Value = 564
$ws.Range("F34").Value = 2144
$ws.Range("F37").Value = 467
$ws.Range("F40").Value = 180
$ws.Range("F41").Value = 345
$ws.Range("F42").Value = 161
$ws.Range("F43").Value = 1
$ws.Range("F44").Value = 21
